$wb = $excel.ActiveWorkbook

# --- Domain_Conversion_Mapper (sheet11): insert "Placement Pixel Size" / Creative_Size row ---
$ws11 = $wb.Worksheets.Item("Domain_Conversion_Mapper")
$ws11.Rows.Item(21).Insert()
$ws11.Range("A21").Value = "Placement Pixel Size"
$ws11.Range("B21").Value = "Creative_Size"
$ws11.Range("C21").Value = $false
$ws11.Range("D21").Value = "VARCHAR"
$ws11.Range("E21").Value = "MATCH"

# --- Domain_Conversion_S3_Mapper (sheet12): insert Creative_Size / Creative_Size row ---
$ws12 = $wb.Worksheets.Item("Domain_Conversion_S3_Mapper")
$ws12.Rows.Item(21).Insert()
$ws12.Range("A21").Value = "Creative_Size"
$ws12.Range("B21").Value = "Creative_Size"
$ws12.Range("C21").Value = $false
$ws12.Range("D21").Value = "VARCHAR"
$ws12.Range("E21").Value = "MATCH"

# --- Domain_Delivery_Mapper (sheet13): insert "Creative Size" / Creative_Size row ---
$ws13 = $wb.Worksheets.Item("Domain_Delivery_Mapper")
$ws13.Rows.Item(21).Insert()
$ws13.Range("A21").Value = "Creative Size"
$ws13.Range("B21").Value = "Creative_Size"
$ws13.Range("C21").Value = $false
$ws13.Range("D21").Value = "VARCHAR"
$ws13.Range("E21").Value = "MATCH"

# --- Domain_Delivery_S3_Mapper (sheet14): insert Creative_Size / Creative_Size row ---
$ws14 = $wb.Worksheets.Item("Domain_Delivery_S3_Mapper")
$ws14.Rows.Item(21).Insert()
$ws14.Range("A21").Value = "Creative_Size"
$ws14.Range("B21").Value = "Creative_Size"
$ws14.Range("C21").Value = $false
$ws14.Range("D21").Value = "VARCHAR"
$ws14.Range("E21").Value = "MATCH"

# Select the new row on each affected sheet (mirrors the row-insert selection left behind in Excel)
$ws11.Range("A21:XFD21").Select()
$ws12.Range("A21:XFD21").Select()
$ws13.Range("A21:XFD21").Select()
$ws14.Range("A21:XFD21").Select()

# Final active sheet/tab matches the workbook's saved state (Domain_Delivery_S3_Mapper active)
$ws14.Activate()
